# Apply the gradrho_files.xlsx edits:
#  - Add A35 / A36 = "nearly" (matching the existing pattern used in A32/A33)
#  - Add five new rows of TODO-style notes in column B (rows 43-46, 48)
#  - Move the active selection/view down to reflect the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rsquo = [char]0x2019

$ws.Range("A35").Value = "nearly"
$ws.Range("A36").Value = "nearly"

$ws.Range("B43").Value = "look at ph_ob file use"
$ws.Range("B44").Value = "Find_Rho should be able to use W" + $rsquo + "s (see get_rho_from_W)"
$ws.Range("B45").Value = "look into find_grad_rho"
$ws.Range("B46").Value = "fix grad_cost_and_rho (see also the test)"
$ws.Range("B48").Value = "farmer_rho_demo.bash"

$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("B49").Select()
